# Updated symbol list on Sun Jan 22 04:53:21 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns with the latest scraped
# values. Each target cell holds plain text (e.g. "300.86", "-0.80%"), so we
# briefly force a Text number format while writing the value to stop Excel's
# COM layer from auto-coercing the string into a number/percentage, then
# restore whatever style the cell had before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "300.86";
    "E2" = "-0.80%";
    "E3" = "6.79%";
    "E4" = "-3.59%";
    "D5" = "0.07740";
    "E5" = "-0.34%";
    "D6" = "2.197";
    "E6" = "-7.12%";
    "E7" = "-0.31%";
    "D8" = "3.989";
    "E8" = "1.09%";
    "D9" = "0.9129";
    "E9" = "-2.08%";
    "D10" = "0.09383";
    "E10" = "-5.09%";
    "D11" = "0.1794";
    "E11" = "-0.02%";
    "D12" = "0.08408";
    "E12" = "-2.33%";
    "D13" = "0.03538";
    "E13" = "6.72%";
    "D14" = "0.09923";
    "E14" = "0.02%";
    "D15" = "0.001474";
    "E15" = "-1.95%";
    "D16" = "0.005723";
    "E16" = "-0.79%";
    "D17" = "3.474";
    "E17" = "0.38%";
    "D18" = "2.053";
    "E18" = "-4.17%";
    "E19" = "2.84%";
    "D20" = "0.1313";
    "E20" = "-1.46%";
    "D21" = "4.556";
    "E21" = "5.77%";
    "D22" = "0.2225";
    "E22" = "-3.26%";
    "D23" = "0.04636";
    "E23" = "1.59%";
    "D24" = "0.001226";
    "E24" = "0.86%";
    "E25" = "1.79%";
    "E26" = "-0.25%";
    "D27" = "0.0004738";
    "E27" = "39.64%";
    "D39" = "0.01750";
    "E39" = "-1.46%";
    "D40" = "0.04677";
    "E40" = "-2.54%";
    "D41" = "0.007800";
    "E41" = "0.69%";
    "E42" = "-1.87%";
    "D43" = "0.007649";
    "E43" = "7.27%";
    "D44" = "0.002285";
    "E44" = "8.86%";
    "D45" = "0.01011";
    "E45" = "10.20%";
    "D46" = "0.00006086";
    "E46" = "-0.45%";
    "D47" = "0.00000000748";
    "E47" = "-0.24%";
    "D48" = "8.660";
    "E48" = "182.76%";
    "E49" = "34.90%";
    "D50" = "0.00002095";
    "E50" = "-0.24%";
    "D51" = "0.0001996";
    "E51" = "-0.24%";
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $originalStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = $originalStyle
}
